# Generate Report for Handback
# This reflects the "handback" of both the zh-cn and de-de localization
# targets: the overview status text flips from "Ready for handoff" to
# "Handed back: in sync with en-US", the handoff datetime that was still
# the zero-date gets a real timestamp, and the per-language "Latest
# Target File" / "Latest Handback File" / "Latest Handback DateTime"
# columns get populated (with the target-file cell becoming a real
# hyperlink, matching the existing "Source File Name" hyperlink style)
# for both rows in the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: status text + handoff datetime + wider status cols
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value2 = "Handed back: in sync with en-US"
$overview.Range("G2:G3").Value2 = "2016-08-23 17:06:38"

# Widen the two status columns that now hold the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.16666666
$overview.Columns.Item(6).ColumnWidth = 29.16666666

# ---------------------------------------------------------------------
# Helper: wire up one language sheet (zh-cn / de-de) with its handback
# info. Positional parameters only -- named parameter binding is not
# reliable against this host.
# ---------------------------------------------------------------------
function Set-HandbackInfo {
    param($sheetName, $handbackDateTime, $xlfRow2, $xlfRow3)

    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C) keeps its text, only widened below.
    $ws.Columns.Item(3).ColumnWidth = 29.16666666

    # Latest Target File (I) / Latest Handback File (J) / Latest Handback
    # DateTime (K) for both data rows.
    $ws.Range("I2").Value2 = "bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.md"
    $ws.Range("J2").Value2 = $xlfRow2
    $ws.Range("K2").Value2 = $handbackDateTime

    $ws.Range("I3").Value2 = "fda9913a-55c6-4958-b398-4e3716a1c1ea.md"
    $ws.Range("J3").Value2 = $xlfRow3
    $ws.Range("K3").Value2 = $handbackDateTime

    # Latest Target File becomes a real hyperlink, just like column A's
    # "Source File Name" hyperlinks (same targets). Rebuild the hyperlink
    # collection in row order (A2, I2, A3, I3) so the new links land right
    # after the matching source-file link for their own row.
    $bdAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14ef90659a43ed09fd0d7f5e511b079f57adfa5/e2e/bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.md"
    $fdaAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c14ef90659a43ed09fd0d7f5e511b079f57adfa5/e2e/fda9913a-55c6-4958-b398-4e3716a1c1ea.md"

    $ws.Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $bdAddress, "", "", "bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $bdAddress, "", "", "bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $fdaAddress, "", "", "fda9913a-55c6-4958-b398-4e3716a1c1ea.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $fdaAddress, "", "", "fda9913a-55c6-4958-b398-4e3716a1c1ea.md")
    # NOTE: Hyperlinks.Add already applies the "Hyperlink" visual style to
    # its target cell, so no extra .Style assignment is needed here (doing
    # so would just strip the underline/colour back off again).

    # Latest Target File / Latest Handback File columns need to be wide
    # enough for the long file names, matching the status column.
    $ws.Columns.Item(9).ColumnWidth = 39.16666666
    $ws.Columns.Item(10).ColumnWidth = 39.16666666
}

# ---------------------------------------------------------------------
# 2. zh-cn sheet
# ---------------------------------------------------------------------
Set-HandbackInfo "zh-cn" "2016-08-23 17:06:38" "bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.dfc3b3f41108d948ac9264ead63eff4781a9e0d0.zh-cn.xlf" "fda9913a-55c6-4958-b398-4e3716a1c1ea.23251e62049597268cc1b201b1760760c0dc13a0.zh-cn.xlf"

# ---------------------------------------------------------------------
# 3. de-de sheet
# ---------------------------------------------------------------------
Set-HandbackInfo "de-de" "2016-08-23 17:06:45" "bd9fa44f-5f28-4ee3-babc-e7aabffa4ade.dfc3b3f41108d948ac9264ead63eff4781a9e0d0.de-de.xlf" "fda9913a-55c6-4958-b398-4e3716a1c1ea.23251e62049597268cc1b201b1760760c0dc13a0.de-de.xlf"

Write-Host "Handback report generated."
